$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new rows: a "# outliers" row after each existing (Bias, RMSE) block ---
# Before: row4 = RJIVE/Bias, row5 = (blank)/RMSE
# Insert at row4 pushes old row4 (RJIVE/Bias) and row5 ((blank)/RMSE) down to row5/row6,
# creating a blank new row4 for the post_lasso "# outliers" line.
$ws.Rows("4:4").Insert()

# Insert at row7 (after the now-shifted RJIVE/RMSE row6) to create a blank new row7 for the
# RJIVE "# outliers" line.
$ws.Rows("7:7").Insert()

# --- Give the new label cell (A8) the same centered look as the other label cells ---
# (copy format only, from the existing "RJIVE" label cell, before writing the new
# "Rliml" text into it)
$ws.Range("A5").Copy()
$ws.Range("A8").PasteSpecial(-4122)

# --- Append a new 3-row block (rows 8-10) for "Rliml" ---
# (written before the "# outliers" rows below so the shared-string table gets
# "Rliml" at index 8 and "# outliers" at index 9, matching the authored file)
$ws.Cells.Item(8,1).Value = "Rliml"
$ws.Cells.Item(8,2).Value = "Bias"
$ws.Cells.Item(8,3).Value = [double]"-1.94882109650196E-3"
$ws.Cells.Item(8,4).Value = [double]"8.4418465313751401E-4"
$ws.Cells.Item(8,5).Value = [double]"-2.65391534145631E-3"
$ws.Cells.Item(8,6).Value = [double]"-1.4844115069809799E-3"

$ws.Cells.Item(9,2).Value = "RMSE"
$ws.Cells.Item(9,3).Value = [double]"6.9427255325198001E-2"
$ws.Cells.Item(9,4).Value = [double]"8.9980141564640806E-2"
$ws.Cells.Item(9,5).Value = [double]"4.3352744458729497E-2"
$ws.Cells.Item(9,6).Value = [double]"4.7017792876855001E-2"

# --- Fill the new "# outliers" row for the post_lasso block (row 4) ---
$ws.Cells.Item(4,2).Value = "# outliers"
$ws.Cells.Item(4,3).Value = 2
$ws.Cells.Item(4,4).Value = 5
$ws.Cells.Item(4,5).Value = 0
$ws.Cells.Item(4,6).Value = 0

# --- Fill the new "# outliers" row for the RJIVE block (row 7) ---
$ws.Cells.Item(7,2).Value = "# outliers"
$ws.Cells.Item(7,3).Value = 0
$ws.Cells.Item(7,4).Value = 0
$ws.Cells.Item(7,5).Value = 0
$ws.Cells.Item(7,6).Value = 0

# --- Fill the new "# outliers" row for the Rliml block (row 10) ---
$ws.Cells.Item(10,2).Value = "# outliers"
$ws.Cells.Item(10,3).Value = 0
$ws.Cells.Item(10,4).Value = 281
$ws.Cells.Item(10,5).Value = 0
$ws.Cells.Item(10,6).Value = 0

# --- Merge the label column for each 3-row block ---
$ws.Range("A2:A4").Merge()
$ws.Range("A5:A7").Merge()
$ws.Range("A8:A10").Merge()

# --- New column width for column H (authored raw width 14.73046875 chars) ---
$ws.Columns("H:H").ColumnWidth = 13.8

# --- Update selection to match the authored state ---
[void]$ws.Range("D10").Select()
